# REPORTGEN-541: align left the "rule name" column in the STIG security
# report summary tables (Cat1/Cat2/Cat3 and Rule 1..8 rows).
#
# The first column of each "Rules"/"STIG V4R8" table inherits a justified
# (both) alignment from the Normal style. We explicitly set these
# paragraphs to left alignment, which emits <w:jc w:val="left"/> in the
# paragraph properties.

$wdAlignParagraphLeft = 0

$d = $word.ActiveDocument

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)

    # Only the summary tables whose header row starts with "STIG V4R8" or
    # "Rules" contain the rule-name column that needs to be realigned.
    $headerText = $t.Cell(1, 1).Range.Text
    $headerText = $headerText -replace "[\x07\x0d]", ""

    if ($headerText -eq "STIG V4R8" -or $headerText -eq "Rules") {
        for ($ri = 2; $ri -le $t.Rows.Count; $ri++) {
            $cell = $t.Cell($ri, 1)
            $cell.Range.ParagraphFormat.Alignment = $wdAlignParagraphLeft
        }
    }
}
